# Atualização de bases das ligas, do dia: 16-06-2024 às 07:16
# Swap the full data (columns B:AD) between each of the following row
# pairs. Column A (the running index) is left untouched on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(38, 39),
    @(110, 111),
    @(237, 238),
    @(249, 250),
    @(268, 269)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$($r1):AD$($r1)")
    $range2 = $ws.Range("B$($r2):AD$($r2)")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
